$wb = $excel.ActiveWorkbook

# Locate the existing "data" sheet so the new sheet is inserted right after it
# (sheetId=2, rId2, second tab) exactly like the target workbook layout.
$dataSheet = $wb.Worksheets.Item("data")

$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# ---- Header row (row 1): bold / bordered / centered header cells ----
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Reuse the header style already used on the "data" sheet (bold, thin border,
# centered/top aligned) so the new header cells share the exact same style
# record instead of creating a near-duplicate one.
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data row (row 2) ----
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Repeat Disorders"
$ws.Range("C2").Value = 3597

# data_version ("0.148") must stay a literal text value, not be coerced into
# the number 0.148 - mark it as Text before writing so it round-trips as a
# string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.148"

$ws.Range("E2").Value = "2021-09-23T03:36:54.755664Z"
$ws.Range("F2").Value = "2021-10-05 14:35:37.044132"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3597/?format=json"

# A2 also carries the bold/bordered style from the header.
$dataSheet.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the Value to A2 since PasteSpecial(Formats) shouldn't touch the
# value, but make sure it stays the number 0 (not blank) after the format copy.
$ws.Range("A2").Value = 0

$ws.Range("A1").Select() | Out-Null
